# Applies the "generateAgesAtDeath adjustment years + defaultInflation" edit
# described in the commit message:
#  - Sheet1 ("Sheet1"): add a "length" column (simulation length), rename
#    "inflation" -> "defaultInflation", and add a "mortadjyears" column for
#    each person (p1mortadjyears, p2mortadjyears).
#  - Sheet2 ("cashflows"): rename "inflationadj" -> "defaultInflationAdj" and
#    add a new "inflation" column (per-cash-flow inflation rate override).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("cashflows")

# ---------------------------------------------------------------------
# Sheet1 - header row (row 1)
# ---------------------------------------------------------------------
$ws1.Cells.Item(1,1).Value  = "description"
$ws1.Cells.Item(1,2).Value  = "nTrials"
$ws1.Cells.Item(1,3).Value  = "startValue"
$ws1.Cells.Item(1,4).Value  = "lengthType"
$ws1.Cells.Item(1,5).Value  = "length"
$ws1.Cells.Item(1,6).Value  = "seed"
$ws1.Cells.Item(1,7).Value  = "defaultInflation"
$ws1.Cells.Item(1,8).Value  = "ror"
$ws1.Cells.Item(1,9).Value  = "stdDev"
$ws1.Cells.Item(1,10).Value = "targetValue"
$ws1.Cells.Item(1,11).Value = "targetValueIsReal"
$ws1.Cells.Item(1,12).Value = "p1name"
$ws1.Cells.Item(1,13).Value = "p1init"
$ws1.Cells.Item(1,14).Value = "p1age"
$ws1.Cells.Item(1,15).Value = "p1gender"
$ws1.Cells.Item(1,16).Value = "p1retage"
$ws1.Cells.Item(1,17).Value = "p1mortfactor"
$ws1.Cells.Item(1,18).Value = "p1mortadjyears"
$ws1.Cells.Item(1,19).Value = "p2name"
$ws1.Cells.Item(1,20).Value = "p2init"
$ws1.Cells.Item(1,21).Value = "p2age"
$ws1.Cells.Item(1,22).Value = "p2gender"
$ws1.Cells.Item(1,23).Value = "p2retage"
$ws1.Cells.Item(1,24).Value = "p2mortfactor"
$ws1.Cells.Item(1,25).Value = "p2mortadjyears"

# ---------------------------------------------------------------------
# Sheet1 - row 2 (sim1)
# ---------------------------------------------------------------------
$ws1.Cells.Item(2,1).Value  = "sim1"
$ws1.Cells.Item(2,2).Value  = 500
$ws1.Cells.Item(2,3).Value  = 1000000
$ws1.Cells.Item(2,4).Value  = "R"
$ws1.Cells.Item(2,5).Value  = 10
$ws1.Cells.Item(2,6).Value  = 101
$ws1.Cells.Item(2,7).Value  = 0.02
$ws1.Cells.Item(2,8).Value  = 0.1
$ws1.Cells.Item(2,9).Value  = 0.08
$ws1.Cells.Item(2,10).Value = 0.0000009999999999999999
$ws1.Cells.Item(2,11).Value = $true
$ws1.Cells.Item(2,12).Value = "Rex"
$ws1.Cells.Item(2,13).Value = "RM"
$ws1.Cells.Item(2,14).Value = 56
$ws1.Cells.Item(2,15).Value = "M"
$ws1.Cells.Item(2,16).Value = 65
$ws1.Cells.Item(2,17).Value = 1
$ws1.Cells.Item(2,18).Value = 0
$ws1.Cells.Item(2,25).Value = 0

# ---------------------------------------------------------------------
# Sheet1 - row 3 (sim2)
# ---------------------------------------------------------------------
$ws1.Cells.Item(3,1).Value  = "sim2"
$ws1.Cells.Item(3,2).Value  = 500
$ws1.Cells.Item(3,3).Value  = 1000000
$ws1.Cells.Item(3,4).Value  = "R"
$ws1.Cells.Item(3,5).Value  = 10
$ws1.Cells.Item(3,6).Value  = 101
$ws1.Cells.Item(3,7).Value  = 0.02
$ws1.Cells.Item(3,8).Value  = 0.1
$ws1.Cells.Item(3,9).Value  = 0.08
$ws1.Cells.Item(3,10).Value = 0.0000009999999999999999
$ws1.Cells.Item(3,11).Value = $true
$ws1.Cells.Item(3,12).Value = "Rex"
$ws1.Cells.Item(3,13).Value = "RM"
$ws1.Cells.Item(3,14).Value = 56
$ws1.Cells.Item(3,15).Value = "M"
$ws1.Cells.Item(3,16).Value = 65
$ws1.Cells.Item(3,17).Value = 1
$ws1.Cells.Item(3,18).Value = 0
$ws1.Cells.Item(3,19).Value = "Julie"
$ws1.Cells.Item(3,20).Value = "JM"
$ws1.Cells.Item(3,21).Value = 53
$ws1.Cells.Item(3,22).Value = "F"
$ws1.Cells.Item(3,23).Value = 65
$ws1.Cells.Item(3,24).Value = 1
$ws1.Cells.Item(3,25).Value = 0

# ---------------------------------------------------------------------
# Sheet2 ("cashflows") - header row (row 1)
# ---------------------------------------------------------------------
$ws2.Cells.Item(1,1).Value  = "simulation"
$ws2.Cells.Item(1,2).Value  = "description"
$ws2.Cells.Item(1,3).Value  = "starttype"
$ws2.Cells.Item(1,4).Value  = "start"
$ws2.Cells.Item(1,5).Value  = "endtype"
$ws2.Cells.Item(1,6).Value  = "end"
$ws2.Cells.Item(1,7).Value  = "type"
$ws2.Cells.Item(1,8).Value  = "amount"
$ws2.Cells.Item(1,9).Value  = "defaultInflationAdj"
$ws2.Cells.Item(1,10).Value = "inflation"

# ---------------------------------------------------------------------
# Sheet2 data rows 2-7: add new "inflation" column J = 0 for every row.
# (all other existing columns are unchanged)
# ---------------------------------------------------------------------
$ws2.Cells.Item(2,10).Value = 0
$ws2.Cells.Item(3,10).Value = 0
$ws2.Cells.Item(4,10).Value = 0
$ws2.Cells.Item(5,10).Value = 0
$ws2.Cells.Item(6,10).Value = 0
$ws2.Cells.Item(7,10).Value = 0

# ---------------------------------------------------------------------
# Restore selections recorded in the saved workbook
# ---------------------------------------------------------------------
$ws1.Range("Y2").Select()
$ws2.Range("J8").Select()
$ws2.Activate()

Write-Host "edit applied"
